$wb = $excel.ActiveWorkbook

# NBR sheet (sheet1) - Reaction_number column C, rows 2-20
$ws1 = $wb.Worksheets.Item("NBR")
$nbrValues = @(663, 647, 637, 641, 635, 638, 634, 628, 584, 582, 584, 581, 573, 570, 564, 562, 0, 551, 547)
for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 3).Value = $nbrValues[$i]
}

# BAR sheet (sheet2) - Reaction_number column C, rows 2-20
$ws2 = $wb.Worksheets.Item("BAR")
$barValues = @(688, 688, 695, 690, 684, 681, 676, 674, 675, 675, 664, 670, 665, 665, 661, 663, 0, 663, 665)
for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 3).Value = $barValues[$i]
}
